$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new track info
$ws.Range("A2").Value = "EL SABOR DE NACHO"
$ws.Range("B2").Value = "QUE SE SEPA"
$ws.Range("C2").Value = "00:03:46"

# Delete rows 3 through 15 (old tracks no longer in the playlist)
$ws.Range("A3:C15").EntireRow.Delete()
